$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (I1, J1) with their corresponding data
# values in row 2 (I2, J2), matching the commit "I0 and IF added".
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header formatting (bold font, thin border, centered)
# by copying the format from the neighboring header cell H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
